$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price / 1h-volume data (and two re-ordered coin
# pairs: Cosmos/Toncoin swapped rows 29-30, ARBITRUM/Celestia swapped
# rows 40-41, EnergySwap/ApeXProtocol swapped rows 45-47) as scraped on
# Thu Jan 25 18:27:43 UTC 2024.
$updates = @(
    @{ Cell = "D2"; Value = "39.964.80" },
    @{ Cell = "E2"; Value = "  +0.70%  " },
    @{ Cell = "D3"; Value = "2.219.78" },
    @{ Cell = "E3"; Value = "  +0.56%  " },
    @{ Cell = "E4"; Value = "  -0.07%  " },
    @{ Cell = "D5"; Value = "292.35" },
    @{ Cell = "E5"; Value = "  +0.48%  " },
    @{ Cell = "D6"; Value = "87.54" },
    @{ Cell = "E6"; Value = "  +2.02%  " },
    @{ Cell = "D7"; Value = "0.512" },
    @{ Cell = "E7"; Value = "  -0.18%  " },
    @{ Cell = "D9"; Value = "0.469" },
    @{ Cell = "E9"; Value = "  -0.02%  " },
    @{ Cell = "D10"; Value = "30.95" },
    @{ Cell = "E10"; Value = "  +0.98%  " },
    @{ Cell = "D11"; Value = "0.0780" },
    @{ Cell = "E11"; Value = "  -0.34%  " },
    @{ Cell = "D12"; Value = "49.97" },
    @{ Cell = "E12"; Value = "  +6.35%  " },
    @{ Cell = "D13"; Value = "0.112" },
    @{ Cell = "E13"; Value = "  +2.98%  " },
    @{ Cell = "D14"; Value = "6.49" },
    @{ Cell = "E14"; Value = "  +2.84%  " },
    @{ Cell = "D15"; Value = "2.554.07" },
    @{ Cell = "E15"; Value = "  -0.61%  " },
    @{ Cell = "D16"; Value = "13.80" },
    @{ Cell = "E16"; Value = "  -1.52%  " },
    @{ Cell = "D17"; Value = "2.250.28" },
    @{ Cell = "E17"; Value = "  +0.91%  " },
    @{ Cell = "D18"; Value = "0.734" },
    @{ Cell = "E18"; Value = "  +1.23%  " },
    @{ Cell = "D19"; Value = "39.900.38" },
    @{ Cell = "E19"; Value = "  +0.38%  " },
    @{ Cell = "D20"; Value = "0.0₃0885" },
    @{ Cell = "E20"; Value = "  +0.73%  " },
    @{ Cell = "D21"; Value = "11.30" },
    @{ Cell = "E21"; Value = "  +2.39%  " },
    @{ Cell = "D22"; Value = "5.77" },
    @{ Cell = "E22"; Value = "  -0.11%  " },
    @{ Cell = "D23"; Value = "65.66" },
    @{ Cell = "E23"; Value = "  +0.52%  " },
    @{ Cell = "D24"; Value = "237.67" },
    @{ Cell = "E24"; Value = "  +1.17%  " },
    @{ Cell = "E25"; Value = "  +0.07%  " },
    @{ Cell = "D26"; Value = "2.46" },
    @{ Cell = "E26"; Value = "  +0.00%  " },
    @{ Cell = "D27"; Value = "1.83" },
    @{ Cell = "E27"; Value = "  +0.27%  " },
    @{ Cell = "D28"; Value = "23.42" },
    @{ Cell = "E28"; Value = "  +3.57%  " },
    @{ Cell = "B29"; Value = "Toncoin" },
    @{ Cell = "C29"; Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton" },
    @{ Cell = "D29"; Value = "2.15" },
    @{ Cell = "E29"; Value = "  -2.52%  " },
    @{ Cell = "B30"; Value = "Cosmos" },
    @{ Cell = "C30"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" },
    @{ Cell = "D30"; Value = "9.25" },
    @{ Cell = "E30"; Value = "  +0.36%  " },
    @{ Cell = "D31"; Value = "156.91" },
    @{ Cell = "E31"; Value = "  +3.78%  " },
    @{ Cell = "D32"; Value = "31.99" },
    @{ Cell = "E32"; Value = "  -2.37%  " },
    @{ Cell = "D33"; Value = "1.00" },
    @{ Cell = "E33"; Value = "  -0.08%  " },
    @{ Cell = "D34"; Value = "4.96" },
    @{ Cell = "E34"; Value = "  +0.85%  " },
    @{ Cell = "D35"; Value = "0.0712" },
    @{ Cell = "E35"; Value = "  -0.36%  " },
    @{ Cell = "D36"; Value = "2.92" },
    @{ Cell = "E36"; Value = "  +4.50%  " },
    @{ Cell = "E37"; Value = "  -1.82%  " },
    @{ Cell = "E38"; Value = "  +0.31%  " },
    @{ Cell = "D39"; Value = "0.0986" },
    @{ Cell = "E39"; Value = "  -0.27%  " },
    @{ Cell = "B40"; Value = "ARBITRUM" },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" },
    @{ Cell = "D40"; Value = "1.71" },
    @{ Cell = "E40"; Value = "  +1.18%  " },
    @{ Cell = "B41"; Value = "Celestia" },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia" },
    @{ Cell = "D41"; Value = "15.40" },
    @{ Cell = "E41"; Value = "  -3.33%  " },
    @{ Cell = "D42"; Value = "2.119.06" },
    @{ Cell = "E42"; Value = "  +3.12%  " },
    @{ Cell = "D43"; Value = "3.75" },
    @{ Cell = "E43"; Value = "  -0.47%  " },
    @{ Cell = "D44"; Value = "0.0270" },
    @{ Cell = "E44"; Value = "  +1.59%  " },
    @{ Cell = "B45"; Value = "EnergySwap" },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" },
    @{ Cell = "D45"; Value = "17.88" },
    @{ Cell = "E45"; Value = "  -0.55%  " },
    @{ Cell = "D46"; Value = "9.96" },
    @{ Cell = "E46"; Value = "  +0.69%  " },
    @{ Cell = "B47"; Value = "ApeXProtocol" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex" },
    @{ Cell = "D47"; Value = "2.09" },
    @{ Cell = "E47"; Value = "  -0.72%  " },
    @{ Cell = "D48"; Value = "2.71" },
    @{ Cell = "E48"; Value = "  +5.27%  " },
    @{ Cell = "D49"; Value = "2.425.75" },
    @{ Cell = "E49"; Value = "  -0.73%  " },
    @{ Cell = "D50"; Value = "1.49" },
    @{ Cell = "E50"; Value = "  +4.36%  " },
    @{ Cell = "D51"; Value = "88.77" },
    @{ Cell = "E51"; Value = "  +0.11%  " }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    # Prefix with an apostrophe so Excel always treats the value as literal
    # text (many of the price figures, e.g. "1.00" or "2.15", would
    # otherwise be auto-converted to numbers).
    $cell.Value = "'" + $update.Value
    # Reset the style to "Normal" so the quote-prefix text entry does not
    # leave behind a stray number-format style on the cell.
    $cell.Style = "Normal"
}
